# Natmi following Dr Hou advice
# Rebuild the LR-pair data rows (2-9) for Rspo3 -> Sdc4.
# The sending cluster "FAPs" keeps its original 4 target-cluster rows
# (now with updated stats), and a new sending cluster "sCs" is added
# with the same 4 target clusters (ECs, FAPs, M2, sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns A..T in order.
$rowsData = @(
    @(2,  "FAPs", "Rspo3", "Sdc4", "ECs", 3, 1, 4.83236,  14.49708,  0.975350813525687, 0.975350813525687, 3, 1, 3.438907666666667,  10.316723,           0.05825422340060618, 0.05825422340060618, 16.61803985209334,  149.56235866884,    0.05681830418508835,  0.05681830418508835),
    @(3,  "FAPs", "Rspo3", "Sdc4", "FAPs", 3, 1, 4.83236,  14.49708,  0.975350813525687, 0.975350813525687, 3, 1, 10.383857,           31.151571,           0.1758999031294962,  0.1758999031294962,  50.17853521252001,  451.60681691268,    0.1715641136164436,   0.1715641136164436),
    @(4,  "FAPs", "Rspo3", "Sdc4", "M2",  3, 1, 4.83236,  14.49708,  0.975350813525687, 0.975350813525687, 3, 1, 25.34077833333333,  76.022335,           0.4292663558501786,  0.4292663558501786,  122.4557635868667,  1102.1018722818,    0.4186852893976787,   0.4186852893976788),
    @(5,  "FAPs", "Rspo3", "Sdc4", "sCs", 3, 1, 4.83236,  14.49708,  0.975350813525687, 0.975350813525687, 3, 1, 19.86921833333334,  59.60765500000001,   0.336579517619719,   0.336579517619719,   96.01521590526669,  864.1369431474002,  0.3282831063264762,   0.3282831063264763),
    @(6,  "sCs",  "Rspo3", "Sdc4", "ECs", 2, 0.6666666666666666, 0.122124, 0.366372, 0.02464918647431296, 0.02464918647431296, 3, 1, 3.438907666666667,  10.316723,           0.05825422340060618, 0.05825422340060618, 0.419973159884,     3.779758438956,     0.001435919215517827, 0.001435919215517828),
    @(7,  "sCs",  "Rspo3", "Sdc4", "FAPs", 2, 0.6666666666666666, 0.122124, 0.366372, 0.02464918647431296, 0.02464918647431296, 3, 1, 10.383857,           31.151571,           0.1758999031294962,  0.1758999031294962,  1.268118152268,     11.413063370412,    0.004335789513052537, 0.004335789513052538),
    @(8,  "sCs",  "Rspo3", "Sdc4", "M2",  2, 0.6666666666666666, 0.122124, 0.366372, 0.02464918647431296, 0.02464918647431296, 3, 1, 25.34077833333333,  76.022335,           0.4292663558501786,  0.4292663558501786,  3.09471721318,      27.85245491862,     0.01058106645249984,  0.01058106645249984),
    @(9,  "sCs",  "Rspo3", "Sdc4", "sCs", 2, 0.6666666666666666, 0.122124, 0.366372, 0.02464918647431296, 0.02464918647431296, 3, 1, 19.86921833333334,  59.60765500000001,   0.336579517619719,   0.336579517619719,   2.42650841974,      21.83857577766,     0.008296411293242758, 0.00829641129324276)
)

foreach ($entry in $rowsData) {
    $r = [int]$entry[0]
    for ($c = 1; $c -lt $entry.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $entry[$c]
    }
}
